$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates driven by the cryptos-list refresh diff.
# Column D holds price strings that LOOK like plain numbers (e.g. "581.53"),
# so Excel would silently coerce a bare .Value assignment into a real number.
# Force the cell to Text format, assign the literal string, then restore the
# default "Normal" style so no stray numFmt/style id is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.859.34'
$ws.Range("E2").Value = '  +0.47%  '
Set-TextValue 'D3' '3.307.78'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.15%  '
Set-TextValue 'D5' '581.53'
$ws.Range("E5").Value = '  -1.89%  '
Set-TextValue 'D6' '175.03'
$ws.Range("E6").Value = '  -6.89%  '
Set-TextValue 'D7' '1.00'
$ws.Range("E7").Value = '  +0.03%  '
Set-TextValue 'D8' '0.580'
$ws.Range("E8").Value = '  -2.38%  '
Set-TextValue 'D9' '3.302.75'
$ws.Range("E9").Value = '  -1.97%  '
Set-TextValue 'D10' '0.174'
$ws.Range("E10").Value = '  -4.76%  '
Set-TextValue 'D11' '0.575'
$ws.Range("E11").Value = '  -2.30%  '
Set-TextValue 'D12' '45.39'
$ws.Range("E12").Value = '  -4.55%  '
Set-TextValue 'D13' '0.0000268'
$ws.Range("E13").Value = '  -2.42%  '
Set-TextValue 'D14' '670.60'
$ws.Range("E14").Value = '  +4.62%  '
Set-TextValue 'D15' '3.846.48'
$ws.Range("E15").Value = '  -1.94%  '
Set-TextValue 'D16' '8.33'
Set-TextValue 'D17' '67.987.70'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("E18").Value = '  -0.75%  '
Set-TextValue 'D19' '3.319.52'
$ws.Range("E19").Value = '  -1.88%  '
Set-TextValue 'D20' '17.40'
$ws.Range("E20").Value = '  -3.55%  '
Set-TextValue 'D21' '10.86'
$ws.Range("E21").Value = '  -2.93%  '
Set-TextValue 'D22' '0.886'
$ws.Range("E22").Value = '  -2.73%  '
Set-TextValue 'D23' '5.41'
$ws.Range("E23").Value = '  +5.69%  '
Set-TextValue 'D24' '17.10'
$ws.Range("E24").Value = '  -5.33%  '
Set-TextValue 'D25' '97.78'
$ws.Range("E25").Value = '  -1.51%  '
Set-TextValue 'D26' '3.86'
$ws.Range("E26").Value = '  -4.16%  '
Set-TextValue 'D27' '2.66'
$ws.Range("E27").Value = '  -6.89%  '
Set-TextValue 'D28' '9.16'
$ws.Range("E28").Value = '  -5.61%  '
Set-TextValue 'D29' '32.88'
$ws.Range("E29").Value = '  +1.31%  '
Set-TextValue 'D30' '8.38'
$ws.Range("E30").Value = '  -3.65%  '
Set-TextValue 'D31' '7.01'
$ws.Range("E31").Value = '  +0.86%  '
Set-TextValue 'D32' '596.43'
$ws.Range("E32").Value = '  -2.43%  '
Set-TextValue 'D33' '10.93'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.103'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D35' '3.747.69'
$ws.Range("E35").Value = '  -6.29%  '
$ws.Range("E36").Value = '  +0.00%  '
Set-TextValue 'D37' '3.39'
$ws.Range("E37").Value = '  -12.71%  '
Set-TextValue 'D38' '55.42'
$ws.Range("E38").Value = '  -1.47%  '
Set-TextValue 'D39' '0.130'
$ws.Range("E39").Value = '  -1.62%  '
Set-TextValue 'D40' '2.61'
$ws.Range("E40").Value = '  -8.26%  '
Set-TextValue 'D41' '32.20'
$ws.Range("E41").Value = '  -4.76%  '
Set-TextValue 'D42' '3.08'
$ws.Range("E42").Value = '  -5.60%  '
Set-TextValue 'D43' '0.0₃0665'
$ws.Range("E43").Value = '  -5.89%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D44' '3.27'
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D45' '0.329'
$ws.Range("E45").Value = '  -4.58%  '
Set-TextValue 'D46' '0.0404'
$ws.Range("E46").Value = '  -4.65%  '
Set-TextValue 'D47' '2.59'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E48").Value = '  -2.05%  '
Set-TextValue 'D49' '1.01'
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("E50").Value = '  -2.47%  '
Set-TextValue 'D51' '2.70'
$ws.Range("E51").Value = '  -4.24%  '
